$d = $word.ActiveDocument
$xml = $d.WordOpenXML

# ---------------------------------------------------------------------------
# 1) Add a new dropdown list item "Religionslehre (jd)" right after the
#    "Religionslehre (orth)" entry and before "Religionslehre (alev)".
# ---------------------------------------------------------------------------
$oldList = '<w:listItem w:displayText="Religionslehre (orth)" w:value="Religionslehre (orth)"/><w:listItem w:displayText="Religionslehre (alev)" w:value="Religionslehre (alev)"/>'
$newList = '<w:listItem w:displayText="Religionslehre (orth)" w:value="Religionslehre (orth)"/><w:listItem w:displayText="Religionslehre (jd)" w:value="Religionslehre (jd)"/><w:listItem w:displayText="Religionslehre (alev)" w:value="Religionslehre (alev)"/>'
if ($xml.IndexOf($oldList) -lt 0) {
    throw "listItem anchor not found"
}
$xml = $xml.Replace($oldList, $newList)

# ---------------------------------------------------------------------------
# 2) Renumber bookmarks 6..9 -> 7..10 (highest id first, so we never clash
#    with a value we already wrote).
# ---------------------------------------------------------------------------
function Rename-Bookmark($xmlText, [string]$oldId, [string]$newId, [string]$bmName) {
    $startOld = '<w:bookmarkStart w:id="' + $oldId + '" w:name="' + $bmName + '"/>'
    $startNew = '<w:bookmarkStart w:id="' + $newId + '" w:name="' + $bmName + '"/>'
    $pos = $xmlText.IndexOf($startOld)
    if ($pos -lt 0) {
        throw ("bookmarkStart not found for " + $bmName)
    }
    $xmlText = $xmlText.Substring(0, $pos) + $startNew + $xmlText.Substring($pos + $startOld.Length)

    $endOld = '<w:bookmarkEnd w:id="' + $oldId + '"/>'
    $endNew = '<w:bookmarkEnd w:id="' + $newId + '"/>'
    $endPos = $xmlText.IndexOf($endOld, $pos)
    if ($endPos -lt 0) {
        throw ("bookmarkEnd not found for " + $bmName)
    }
    $xmlText = $xmlText.Substring(0, $endPos) + $endNew + $xmlText.Substring($endPos + $endOld.Length)

    return $xmlText
}

$xml = Rename-Bookmark $xml "9" "10" "Text13"
$xml = Rename-Bookmark $xml "8" "9"  "Text12"
$xml = Rename-Bookmark $xml "7" "8"  "Text14"
$xml = Rename-Bookmark $xml "6" "7"  "Text7"

# ---------------------------------------------------------------------------
# 3) Drop the stray "_GoBack" bookmark that used to sit inside the Text13
#    field (it collided with id 9/10 and is no longer needed there).
# ---------------------------------------------------------------------------
$oldGoBack = '<w:bookmarkStart w:id="10" w:name="_GoBack"/><w:bookmarkEnd w:id="10"/>'
if ($xml.IndexOf($oldGoBack) -lt 0) {
    throw "old _GoBack bookmark not found"
}
$xml = $xml.Replace($oldGoBack, '')

# ---------------------------------------------------------------------------
# 4) Re-create the "_GoBack" bookmark (now id 6) right after the Text6
#    bookmark end, reflecting where the cursor was when the document was
#    last saved.
# ---------------------------------------------------------------------------
$anchorOld = '<w:bookmarkEnd w:id="5"/>'
$anchorNew = '<w:bookmarkStart w:id="6" w:name="_GoBack"/><w:bookmarkEnd w:id="5"/><w:bookmarkEnd w:id="6"/>'
if ($xml.IndexOf($anchorOld) -lt 0) {
    throw "bookmarkEnd id 5 not found"
}
$xml = $xml.Replace($anchorOld, $anchorNew)

# ---------------------------------------------------------------------------
# 5) Register the new editing session rsid in the glossary part's settings.
# ---------------------------------------------------------------------------
$rsidOld = '<w:rsid w:val="00D95998"/><w:rsid w:val="00ED42D6"/>'
$rsidNew = '<w:rsid w:val="00D95998"/><w:rsid w:val="00E94EB3"/><w:rsid w:val="00ED42D6"/>'
if ($xml.IndexOf($rsidOld) -lt 0) {
    throw "glossary rsid anchor not found"
}
$xml = $xml.Replace($rsidOld, $rsidNew)

$d.WordOpenXML = $xml
